$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Structural inserts -------------------------------------------------
# Insert a new row 6 (pushes old row7->8, row8->9, row9->10)
$ws.Rows.Item(6).Insert()

# Insert two new columns at K (pushes old K,L onward two to the right)
$ws.Columns.Item(11).Insert()
$ws.Columns.Item(11).Insert()

Write-Output "structure done"
